# Edit: add a "Player Info" sheet (as the first tab) with a single player's
# bio data, and on the existing "ODI Batting" / "ODI Bowling" sheets replace
# the MATCH_CARD_LINK column (a full scorecard URL) with a terser MATCH_CODE
# column (just the numeric code parsed out of that URL).

$wb = $excel.ActiveWorkbook

function Convert-MatchLinkSheet($ws) {
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count

    $linkCol = 0
    for ($c = 1; $c -le $ncols; $c++) {
        $header = $ws.Cells.Item(1, $c).Value2
        if ($header -eq "MATCH_CARD_LINK") {
            $linkCol = $c
        }
    }
    if ($linkCol -eq 0) {
        return
    }

    $ws.Cells.Item(1, $linkCol).Value = "MATCH_CODE"

    for ($r = 2; $r -le $nrows; $r++) {
        $cell = $ws.Cells.Item($r, $linkCol)
        $val = $cell.Value2
        if ($val -match "MatchCode=(\d+)") {
            $code = $matches[1]
            # Force the numeric-looking code to be stored as text (matching
            # the other plain data columns, which carry no cell style) rather
            # than letting it be auto-coerced to a Number.
            $cell.NumberFormat = "@"
            $cell.Value = $code
            $cell.Style = "Normal"
        }
    }
}

Convert-MatchLinkSheet $wb.Worksheets.Item("ODI Batting")
Convert-MatchLinkSheet $wb.Worksheets.Item("ODI Bowling")

# Insert the new "Player Info" sheet as the first tab.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "3724"
$idCell.Style = "Normal"

$playerInfo.Range("B2").Value = "Joseph L Denly"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"
